$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style of an existing header cell (H1) to the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Data values for columns I (I0) and J (IF)
$dataI = @(5, 6, 9, 6, 9, 7, 4, 7, 4, 8, 6, 6, 5)
$dataJ = @(6, 6, 9, 7, 9, 8, 5, 7, 4, 8, 6, 6, 5)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
